{"js": "// Title paragraph: update the date line.\nconst body = context.document.body;\nconst titlePara = body.paragraphs.getFirst();\ntitlePara.load(\"text\");\nawait context.sync();\ntitlePara.insertText('2025-11-06 Thursday', Word.InsertLocation.replace);\n\n// Table of division problems: update each cell value in place,\n// preserving the existing run formatting.\nconst table = body.tables.getFirst();\nconst cellUpdates = [\n  { row: 0, col: 0, value: '69\u00f75=' },\n  { row: 0, col: 1, value: '72\u00f74=' },\n  { row: 0, col: 2, value: '66\u00f78=' },\n  { row: 0, col: 3, value: '27\u00f77=' },\n  { row: 0, col: 4, value: '91\u00f75=' },\n  { row: 4, col: 0, value: '23\u00f75=' },\n  { row: 4, col: 1, value: '90\u00f76=' },\n  { row: 4, col: 2, value: '25\u00f78=' },\n  { row: 4, col: 3, value: '20\u00f78=' },\n  { row: 4, col: 4, value: '81\u00f76=' },\n  { row: 8, col: 0, value: '82\u00f77=' },\n  { row: 8, col: 1, value: '56\u00f77=' },\n  { row: 8, col: 2, value: '66\u00f73=' },\n  { row: 8, col: 3, value: '71\u00f75=' },\n  { row: 8, col: 4, value: '44\u00f79=' },\n  { row: 12, col: 0, value: '38\u00f77=' },\n  { row: 12, col: 1, value: '41\u00f75=' },\n  { row: 12, col: 2, value: '82\u00f78=' },\n  { row: 12, col: 3, value: '36\u00f79=' },\n  { row: 12, col: 4, value: '32\u00f72=' },\n  { row: 16, col: 0, value: '30\u00f77=' },\n  { row: 16, col: 1, value: '82\u00f77=' },\n  { row: 16, col: 2, value: '23\u00f76=' },\n  { row: 16, col: 3, value: '49\u00f77=' },\n  { row: 16, col: 4, value: '22\u00f76=' },\n];\nfor (const u of cellUpdates) {\n  table.getCell(u.row, u.col).value = u.value;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Title paragraph: update the date line.\n$d.Paragraphs.Item(1).Range.Text = \"2025-11-06 Thursday\"\n\n# Table of division problems: update each cell value in place,\n# preserving the existing run formatting (COM is 1-based).\n$table = $d.Tables.Item(1)\n$cellUpdates = @(\n    @{ Row = 1; Col = 1; Value = \"69\u00f75=\" }\n    @{ Row = 1; Col = 2; Value = \"72\u00f74=\" }\n    @{ Row = 1; Col = 3; Value = \"66\u00f78=\" }\n    @{ Row = 1; Col = 4; Value = \"27\u00f77=\" }\n    @{ Row = 1; Col = 5; Value = \"91\u00f75=\" }\n    @{ Row = 5; Col = 1; Value = \"23\u00f75=\" }\n    @{ Row = 5; Col = 2; Value = \"90\u00f76=\" }\n    @{ Row = 5; Col = 3; Value = \"25\u00f78=\" }\n    @{ Row = 5; Col = 4; Value = \"20\u00f78=\" }\n    @{ Row = 5; Col = 5; Value = \"81\u00f76=\" }\n    @{ Row = 9; Col = 1; Value = \"82\u00f77=\" }\n    @{ Row = 9; Col = 2; Value = \"56\u00f77=\" }\n    @{ Row = 9; Col = 3; Value = \"66\u00f73=\" }\n    @{ Row = 9; Col = 4; Value = \"71\u00f75=\" }\n    @{ Row = 9; Col = 5; Value = \"44\u00f79=\" }\n    @{ Row = 13; Col = 1; Value = \"38\u00f77=\" }\n    @{ Row = 13; Col = 2; Value = \"41\u00f75=\" }\n    @{ Row = 13; Col = 3; Value = \"82\u00f78=\" }\n    @{ Row = 13; Col = 4; Value = \"36\u00f79=\" }\n    @{ Row = 13; Col = 5; Value = \"32\u00f72=\" }\n    @{ Row = 17; Col = 1; Value = \"30\u00f77=\" }\n    @{ Row = 17; Col = 2; Value = \"82\u00f77=\" }\n    @{ Row = 17; Col = 3; Value = \"23\u00f76=\" }\n    @{ Row = 17; Col = 4; Value = \"49\u00f77=\" }\n    @{ Row = 17; Col = 5; Value = \"22\u00f76=\" }\n)\n\nforeach ($u in $cellUpdates) {\n    $table.Cell($u.Row, $u.Col).Range.Text = $u.Value\n}\n\n"}
